$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.16"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'24.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'5.047"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.05621"
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Value = "'2.984"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Value = "'0.8393"
$ws.Range("D9").Style = "Normal"
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = "'0.1336"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '9WazirXWRX'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = "'0.06942"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '10MandalaExchangeTokenMDX'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = "'0.02821"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '11BitrueCoinBTR'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = "'0.09412"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '12BitMartTokenBMX'
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = "'0.001508"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '13BitForexTokenBF'
$ws.Range("B15").Value = 'One'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D15").Value = "'0.0005987"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '14OneONE'
$ws.Range("D16").Value = "'0.006270"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'3.499"
$ws.Range("D17").Style = "Normal"
$ws.Range("D20").Value = "'0.03299"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Value = "'3.736"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'0.04673"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'0.1369"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'0.001241"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'0.004527"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'0.00009693"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '26NitroExNTXBestin24h'
$ws.Range("D40").Value = "'0.03630"
$ws.Range("D40").Style = "Normal"
$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D41").Value = "'0.1052"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '40BKEXTokenBKK'
$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D42").Value = "'0.002721"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '41CEJICEJI'
$ws.Range("B43").Value = 'KickToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D43").Value = "'0.003366"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '42KickTokenKICK'
$ws.Range("D44").Value = "'0.008362"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005263"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("D46").Style = "Normal"
